# Update column C ("Förändrad") date serial value from 45202 to 45203
# for rows 2 through 11 on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45202) {
        $cell.Value2 = 45203
    }
}
